# Update the "取得日時" (retrieved datetime) timestamps in column A
# for rows 2-18 on the active sheet ("ランサーズ") from
# "2025-09-27 01:13:25" to "2025-09-27 01:38:04".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-27 01:38:04"

for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
